$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AuthorList")

# Remove the rows for: J.V. Carrion, J. Generowicz, J. Munoz Vidal, J. Rodriguez
# Rows must be deleted from bottom to top so earlier row numbers stay valid.
$ws.Rows.Item(91).Delete()
$ws.Rows.Item(75).Delete()
$ws.Rows.Item(43).Delete()
$ws.Rows.Item(21).Delete()

# Reflect the scrolled/selected view state after the edit.
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A73:XFD73").Select()
